$d = $word.ActiveDocument

# --- Insert Author / Date / AbstractTitle / Abstract paragraphs after the Title paragraph (paragraph 1) ---
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$p = $d.Paragraphs(2)
$p.Range.Text = "Eduardo Santana"
$p.Style = "Author"

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(3)
$p.Range.Text = "2025-09-01"
$p.Style = "Date"

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(4)
$p.Range.Text = "Abstract"
$p.Style = "AbstractTitle"

$p.Range.InsertParagraphAfter()
$p = $d.Paragraphs(5)
$p.Range.Text = "This report contains my self-research related to conscienciology volunteer work."
$p.Style = "Abstract"

# --- Update the Heading2 "Section" text to "1 Section" ---
$d.Content.Find.Execute("Section", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1 Section", 2)
